$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) After the "Spawn Entity" user-story paragraph (the one ending in
#    "...removed first."), add five new top-level list items describing the
#    game's features: Start game / Move character / Shoot enemies /
#    Earn points / Load mods.
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(15)
$newItems = @("Start game", "Move character", "Shoot enemies", "Earn points", "Load mods")

foreach ($itemText in $newItems) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $anchor.Next()
    $anchor.Range.Text = $itemText
    $anchor.Range.ListFormat.ListLevelNumber = 1
}

# ---------------------------------------------------------------------------
# 2) Under "Application Prototype" (right after the blank paragraph that
#    follows the heading), add a new list item "Open mod app" before the
#    "Create Entity" item, carrying Word's usual "_GoBack" last-edit
#    bookmark.
# ---------------------------------------------------------------------------
# $anchor currently sits on "Load mods" -- walk forward to the heading and
# then to the blank paragraph right after it.
$headingPara = $anchor.Next()
$blankAfterHeading = $headingPara.Next()

$blankAfterHeading.Range.InsertParagraphAfter()
$openModApp = $blankAfterHeading.Next()
$openModApp.Range.ListFormat.ListLevelNumber = 1
$openModApp.Range.Text = "Open mod app"
$bm = $d.Range($openModApp.Range.Start, $openModApp.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bm)
